$d = $word.ActiveDocument

# Change 1: merge " " + "восстанавливающегося" (with proofErr wrapping removed)
# into a single run " восстанавливающегося"
$d.Content.Find.Execute(" восстанавливающегося", $true, $false, $false, $false, $false,
                         $true, 1, $false, " восстанавливающегося", 2)

# Change 2: merge the split runs describing the postcondition sentence into one run
$d.Content.Find.Execute("Запрос клиента на определённый ресурс остаётся в очеред",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Запрос клиента на определённый ресурс остаётся в очеред", 2)

$d.Content.Find.Execute("и ожидания, время ожидания запроса обнуляется.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "и ожидания, время ожидания запроса обнуляется.", 2)
